# Target: sheet "汽車" (Car) — add metadata columns (H:N) to match the other
# property sheets (存款/股票/保險) and replace the stray duplicate-data row 1
# with a proper header row, per commit "#5: property boat&car done".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1: header row (B1:N1) -------------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Match the header formatting already used by B1:G1 (bold, centered, bordered)
$ws.Range("H1:N1").Font.Bold = $true
$ws.Range("H1:N1").HorizontalAlignment = -4108
$ws.Range("H1:N1").VerticalAlignment = -4160
$ws.Range("H1:N1").Borders.LineStyle = 1

# J2/J3 hold a literal "2013-12-12" text (like J column in the other
# property sheets) — force text format first so COM doesn't coerce it to a
# date serial number.
$ws.Range("J2:J3").NumberFormat = "@"

# --- Row 2 -----------------------------------------------------------------
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2013-12-12"
$ws.Range("K2").Value = "王進士"
$ws.Range("L2").Value = 1701
$ws.Range("M2").Value = "tmp93a21"
$ws.Range("N2").Value = 32

# --- Row 3 -----------------------------------------------------------------
$ws.Range("H3").Value = "land"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2013-12-12"
$ws.Range("K3").Value = "王進士"
$ws.Range("L3").Value = 1701
$ws.Range("M3").Value = "tmp93a21"
$ws.Range("N3").Value = 33
